$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$ws1 = $wb.Worksheets.Item("Schedule")
$ws1.Range("E2").Value = 1739.1522525
$ws1.Range("F2").Value = 28.75582428075397

# --- Sheet "Detailed" ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("B21").Value = 8.18614
$ws2.Range("B22").Value = 10.3824
$ws2.Range("C23").Value = "historical"
$ws2.Range("B24").Value = 0.51
$ws2.Range("C24").Value = "historical"
$ws2.Range("C25").Value = "historical"
$ws2.Range("B26").Value = 55.43646
$ws2.Range("C26").Value = "historical"
$ws2.Range("B27").Value = 51.22117
$ws2.Range("B28").Value = 36.06009
$ws2.Range("B29").Value = 36.0601
$ws2.Range("B30").Value = 44.79373
$ws2.Range("B31").Value = 56.98
$ws2.Range("B32").Value = 39.93588
$ws2.Range("B33").Value = 40.20882
$ws2.Range("B34").Value = 43.87656
$ws2.Range("B35").Value = 44.10049
$ws2.Range("B36").Value = 49.23158
$ws2.Range("B37").Value = 23.49459
$ws2.Range("B38").Value = 66.7022
$ws2.Range("B39").Value = 63.34815
$ws2.Range("B40").Value = 106.89055
$ws2.Range("B41").Value = 120.01
$ws2.Range("B42").Value = 125.00189
$ws2.Range("B43").Value = 101.25
$ws2.Range("B44").Value = 85.08304
$ws2.Range("B45").Value = 80.02
$ws2.Range("B46").Value = 69.09674
$ws2.Range("B47").Value = 57.09
$ws2.Range("B48").Value = 57.09
$ws2.Range("B49").Value = 57.09
